$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# CNV / Aneuploidy calls table (row 1 is the header).
# Rows 2-4 are "Aneuploidy" calls whose Event Start / Event End / Event
# Size columns were left blank and whose Control Fractional Copy Number
# column was also left blank - fill those with the literal placeholders
# used elsewhere in the report ("<NA>" for the event columns, "NA" for
# the control-copy-number column). Rows 5-14 are "CNV" calls: their
# Event Start/End/Size numbers carried a spurious ".000" suffix that is
# trimmed off, and their Control Fractional Copy Number column (also
# blank) is filled with "NA" as well.

$t.Cell(2, 7).Range.Text = "<NA>"
$t.Cell(2, 8).Range.Text = "<NA>"
$t.Cell(2, 9).Range.Text = "<NA>"
$t.Cell(2, 11).Range.Text = "NA"
$t.Cell(3, 7).Range.Text = "<NA>"
$t.Cell(3, 8).Range.Text = "<NA>"
$t.Cell(3, 9).Range.Text = "<NA>"
$t.Cell(3, 11).Range.Text = "NA"
$t.Cell(4, 7).Range.Text = "<NA>"
$t.Cell(4, 8).Range.Text = "<NA>"
$t.Cell(4, 9).Range.Text = "<NA>"
$t.Cell(4, 11).Range.Text = "NA"
$t.Cell(5, 7).Range.Text = "18514"
$t.Cell(5, 8).Range.Text = "39512270"
$t.Cell(5, 9).Range.Text = "39493756"
$t.Cell(5, 11).Range.Text = "NA"
$t.Cell(6, 7).Range.Text = "42178604"
$t.Cell(6, 8).Range.Text = "46404297"
$t.Cell(6, 9).Range.Text = "4225693"
$t.Cell(6, 11).Range.Text = "NA"
$t.Cell(7, 7).Range.Text = "47686678"
$t.Cell(7, 8).Range.Text = "98844380"
$t.Cell(7, 9).Range.Text = "51157702"
$t.Cell(7, 11).Range.Text = "NA"
$t.Cell(8, 7).Range.Text = "101275106"
$t.Cell(8, 8).Range.Text = "132148913"
$t.Cell(8, 9).Range.Text = "30873807"
$t.Cell(8, 11).Range.Text = "NA"
$t.Cell(9, 7).Range.Text = "32821561"
$t.Cell(9, 8).Range.Text = "61861320"
$t.Cell(9, 9).Range.Text = "29039759"
$t.Cell(9, 11).Range.Text = "NA"
$t.Cell(10, 7).Range.Text = "12406577"
$t.Cell(10, 8).Range.Text = "43121077"
$t.Cell(10, 9).Range.Text = "30714500"
$t.Cell(10, 11).Range.Text = "NA"
$t.Cell(11, 7).Range.Text = "19314"
$t.Cell(11, 8).Range.Text = "46400789"
$t.Cell(11, 9).Range.Text = "46381475"
$t.Cell(11, 11).Range.Text = "NA"
$t.Cell(12, 7).Range.Text = "49709239"
$t.Cell(12, 8).Range.Text = "70107744"
$t.Cell(12, 9).Range.Text = "20398505"
$t.Cell(12, 11).Range.Text = "NA"
$t.Cell(13, 7).Range.Text = "71001016"
$t.Cell(13, 8).Range.Text = "181472714"
$t.Cell(13, 9).Range.Text = "110471698"
$t.Cell(13, 11).Range.Text = "NA"
$t.Cell(14, 7).Range.Text = "74488138"
$t.Cell(14, 8).Range.Text = "126892338"
$t.Cell(14, 9).Range.Text = "52404200"
$t.Cell(14, 11).Range.Text = "NA"
